# Automatische test-sync: 2025-06-30 20:03:50
#
# Adds Testmail #12 to the "Logs" sheet (row 12), extends the conditional
# formatting ranges to cover the new row, reorders/extends the "Dashboard"
# pivot-style summary table (swap rows 4/5, add a new "Klacht / Probleem"
# row 7) and extends the bar chart's category/value series ranges to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 12 with the new test-mail entry.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Ik wacht nog steeds op reactie. Wanneer hoor ik iets?"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Testmail #12: Ik wacht nog steeds op reactie. Wanneer hoor ik iets?"
$logs.Range("D12").Value = "Klacht / Probleem"
$logs.Range("E12").Value = "Beste klant,`nDank u wel voor uw e-mail. Excuses voor het ongemak dat u heeft ervaren. Om u beter van dienst te kunnen zijn, kunt u alstublieft uw gebruikersnaam vermelden zodat we het specifieke probleem kunnen onderzoeken en een passende oplossing kunnen bieden.`nWij streven ernaar om binnen 24 uur op al onze e-mails te reageren. Zodra we meer informatie hebben, zullen we direct contact met u opnemen.`nMet vriendelijke groet,`n[Naam van het bedrijf] E-mailassistent"
$logs.Range("F12").Value = "2025-06-30 20:03:29"
$logs.Range("G12").Value = "Ja"
$logs.Range("H12").Value = "Nee"
$logs.Range("I12").Value = "Ja"
$logs.Range("J12").Value = "Nee"

# Extend the conditional-formatting sqref ranges (D/G/H/I/J, rows 2:11 ->
# 2:12) so the newly appended row picks up the same highlighting rules.
$logs.Range("D2:D11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D12"))
$logs.Range("G2:G11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G12"))
$logs.Range("H2:H11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H12"))
$logs.Range("I2:I11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I12"))
$logs.Range("J2:J11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J12"))

# ---------------------------------------------------------------------
# 2) Dashboard sheet: swap rows 4/5 ("Openingstijden / Locatie" and
#    "Bestelling / Levering" trade places) and append a new row 7 for
#    "Klacht / Probleem".
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("A5").Value = "Openingstijden / Locatie"

$dash.Range("A7").Value = "Klacht / Probleem"
$dash.Range("B7").Value = 1

# ---------------------------------------------------------------------
# 3) Chart: extend the bar chart's category/value series so it also
#    plots the new Dashboard row (A2:A6/B2:B6 -> A2:A7/B2:B7).
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$7,'Dashboard'!`$B`$2:`$B`$7,1)"
